$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 25, pushing the existing rows 25-63 down to 26-64
# (this also carries the date-format style of column D along for the ride).
$ws.Rows("25:25").Insert()

# Populate the newly inserted row 25 with the latest weekly price record.
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 44645
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = 100114002
$ws.Range("G25").Value = "Camote"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 30
$ws.Range("K25").Value = 18000
$ws.Range("L25").Value = 18000
$ws.Range("M25").Value = 18000
$ws.Range("N25").Value = "$/caja 15 kilos granel"
$ws.Range("O25").Value = "Perú"
$ws.Range("P25").Value = 1200
$ws.Range("Q25").Value = 15
$ws.Range("R25").Value = "Hortaliza"
